$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.110.90'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.256.03'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '397.64'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.52'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.97%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.618'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.19'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0951'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.81%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.772.63'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.23'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '18.93'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.250.48'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.98'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '56.911.40'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.29'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000107'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.90'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '294.23'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.84'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.17'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '28.03'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.38'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.85'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.40'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.168'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.82%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.17'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '40.00'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +10.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0492'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.26'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.45'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '136.81'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.49%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.283'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.86'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.89'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.66'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.29'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.21'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.16%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.46'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.137.54'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -8.22%  '
